$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.499283333333334
$ws.Range("H2").Value = 13.49785
$ws.Range("I2").Value = 0.583046232715136
$ws.Range("J2").Value = 0.6111950921276581
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 137.0763058517778
$ws.Range("R2").Value = 1233.686752666
$ws.Range("S2").Value = 0.1274177863208488
$ws.Range("T2").Value = 0.1425022149262821

$ws.Range("G3").Value = 4.499283333333334
$ws.Range("H3").Value = 13.49785
$ws.Range("I3").Value = 0.583046232715136
$ws.Range("J3").Value = 0.6111950921276581
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("Q3").Value = 156.2967163358666
$ws.Range("R3").Value = 1406.6704470228
$ws.Range("S3").Value = 0.1452839094326636
$ws.Range("T3").Value = 0.1624834293947888

$ws.Range("G4").Value = 4.499283333333334
$ws.Range("H4").Value = 13.49785
$ws.Range("I4").Value = 0.583046232715136
$ws.Range("J4").Value = 0.6111950921276581
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 106.5897018991111
$ws.Range("R4").Value = 959.3073170920001
$ws.Range("S4").Value = 0.09907929584321939
$ws.Range("T4").Value = 0.1108088558016715

$ws.Range("G5").Value = 4.499283333333334
$ws.Range("H5").Value = 13.49785
$ws.Range("I5").Value = 0.583046232715136
$ws.Range("J5").Value = 0.6111950921276581
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 117.957110495675
$ws.Range("R5").Value = 707.7426629740501
$ws.Range("S5").Value = 0.1096457466282656
$ws.Range("T5").Value = 0.08175081466480824

$ws.Range("G6").Value = 4.499283333333334
$ws.Range("H6").Value = 13.49785
$ws.Range("I6").Value = 0.583046232715136
$ws.Range("J6").Value = 0.6111950921276581
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 109.3224526139333
$ws.Range("R6").Value = 983.9020735253999
$ws.Range("S6").Value = 0.1016194944901385
$ws.Range("T6").Value = 0.1136497773401074

$ws.Range("I7").Value = 0.2184296534786591
$ws.Range("J7").Value = 0.228975207608499
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 51.35361195613334
$ws.Range("R7").Value = 462.1825076052
$ws.Range("S7").Value = 0.04773519036984992
$ws.Range("T7").Value = 0.05338634859424097

$ws.Range("I8").Value = 0.2184296534786591
$ws.Range("J8").Value = 0.228975207608499
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("Q8").Value = 58.55425466024
$ws.Range("S8").Value = 0.05442846932673058
$ws.Range("T8").Value = 0.06087201526618582

$ws.Range("I9").Value = 0.2184296534786591
$ws.Range("J9").Value = 0.228975207608499
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 39.93225638693335
$ws.Range("R9").Value = 359.3903074824
$ws.Range("S9").Value = 0.03711859376427479
$ws.Range("T9").Value = 0.04151290003609601

$ws.Range("I10").Value = 0.2184296534786591
$ws.Range("J10").Value = 0.228975207608499
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 44.190888003735
$ws.Range("R10").Value = 265.14532802241
$ws.Range("S10").Value = 0.04107715837540162
$ws.Range("T10").Value = 0.03062673441122567

$ws.Range("I11").Value = 0.2184296534786591
$ws.Range("J11").Value = 0.228975207608499
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 40.95604105132
$ws.Range("R11").Value = 368.60436946188
$ws.Range("S11").Value = 0.03807024164240215
$ws.Range("T11").Value = 0.04257720930075051

$ws.Range("G12").Value = 0.269395
$ws.Range("H12").Value = 0.808185
$ws.Range("I12").Value = 0.03490994636826474
$ws.Range("J12").Value = 0.03659536189327867
$ws.Range("M12").Value = 30.46625333333334
$ws.Range("N12").Value = 91.39876000000001
$ws.Range("O12").Value = 0.2185380492512374
$ws.Range("P12").Value = 0.2331534018544084
$ws.Range("Q12").Value = 8.207456316733335
$ws.Range("R12").Value = 73.86710685060001
$ws.Range("S12").Value = 0.007629151578785896
$ws.Range("T12").Value = 0.008532333117511106

$ws.Range("G13").Value = 0.269395
$ws.Range("H13").Value = 0.808185
$ws.Range("I13").Value = 0.03490994636826474
$ws.Range("J13").Value = 0.03659536189327867
$ws.Range("O13").Value = 0.2491807703757967
$ws.Range("P13").Value = 0.2658454419670822
$ws.Range("Q13").Value = 9.358280147719999
$ws.Range("R13").Value = 84.22452132948
$ws.Range("S13").Value = 0.008698887329821952
$ws.Range("T13").Value = 0.009728710156463986

$ws.Range("G14").Value = 0.269395
$ws.Range("H14").Value = 0.808185
$ws.Range("I14").Value = 0.03490994636826474
$ws.Range("J14").Value = 0.03659536189327867
$ws.Range("M14").Value = 23.69037333333334
$ws.Range("N14").Value = 71.07112000000001
$ws.Range("O14").Value = 0.1699338582153697
$ws.Range("P14").Value = 0.181298667526812
$ws.Range("Q14").Value = 6.382068124133334
$ws.Range("R14").Value = 57.43861311720001
$ws.Range("S14").Value = 0.005932381876450861
$ws.Range("T14").Value = 0.006634690348912895

$ws.Range("G15").Value = 0.269395
$ws.Range("H15").Value = 0.808185
$ws.Range("I15").Value = 0.03490994636826474
$ws.Range("J15").Value = 0.03659536189327867
$ws.Range("M15").Value = 26.2168665
$ws.Range("N15").Value = 52.433733
$ws.Range("O15").Value = 0.18805669340777
$ws.Range("P15").Value = 0.1337556791894743
$ws.Range("Q15").Value = 7.062692750767501
$ws.Range("R15").Value = 42.37615650460501
$ws.Range("S15").Value = 0.006565049081058454
$ws.Range("T15").Value = 0.004894837485220095

$ws.Range("G16").Value = 0.269395
$ws.Range("H16").Value = 0.808185
$ws.Range("I16").Value = 0.03490994636826474
$ws.Range("J16").Value = 0.03659536189327867
$ws.Range("M16").Value = 24.297748
$ws.Range("N16").Value = 72.893244
$ws.Range("O16").Value = 0.1742906287498262
$ws.Range("P16").Value = 0.1859468094622229
$ws.Range("Q16").Value = 6.545691822459999
$ws.Range("R16").Value = 58.91122640214
$ws.Range("S16").Value = 0.006084476502147573
$ws.Range("T16").Value = 0.006804790785170583

$ws.Range("G17").Value = 1.0662095
$ws.Range("H17").Value = 2.132419
$ws.Range("I17").Value = 0.138166322546203
$ws.Range("J17").Value = 0.0965578982697073
$ws.Range("M17").Value = 30.46625333333334
$ws.Range("N17").Value = 91.39876000000001
$ws.Range("O17").Value = 0.2185380492512374
$ws.Range("P17").Value = 0.2331534018544084
$ws.Range("Q17").Value = 32.48340873340667
$ws.Range("R17").Value = 194.90045240044
$ws.Range("S17").Value = 0.03019459860146447
$ws.Range("T17").Value = 0.02251280245749416

$ws.Range("G18").Value = 1.0662095
$ws.Range("H18").Value = 2.132419
$ws.Range("I18").Value = 0.138166322546203
$ws.Range("J18").Value = 0.0965578982697073
$ws.Range("O18").Value = 0.2491807703757967
$ws.Range("P18").Value = 0.2658454419670822
$ws.Range("Q18").Value = 37.038130615492
$ws.Range("R18").Value = 222.228783692952
$ws.Range("S18").Value = 0.03442839069205367
$ws.Range("T18").Value = 0.0256694771409229

$ws.Range("G19").Value = 1.0662095
$ws.Range("H19").Value = 2.132419
$ws.Range("I19").Value = 0.138166322546203
$ws.Range("J19").Value = 0.0965578982697073
$ws.Range("M19").Value = 23.69037333333334
$ws.Range("N19").Value = 71.07112000000001
$ws.Range("O19").Value = 0.1699338582153697
$ws.Range("P19").Value = 0.181298667526812
$ws.Range("Q19").Value = 25.25890110654667
$ws.Range("R19").Value = 151.55340663928
$ws.Range("S19").Value = 0.02347913626570551
$ws.Range("T19").Value = 0.0175058182954874

$ws.Range("G20").Value = 1.0662095
$ws.Range("H20").Value = 2.132419
$ws.Range("I20").Value = 0.138166322546203
$ws.Range("J20").Value = 0.0965578982697073
$ws.Range("M20").Value = 26.2168665
$ws.Range("N20").Value = 52.433733
$ws.Range("O20").Value = 0.18805669340777
$ws.Range("P20").Value = 0.1337556791894743
$ws.Range("Q20").Value = 27.95267212253175
$ws.Range("R20").Value = 111.810688490127
$ws.Range("S20").Value = 0.02598310175835036
$ws.Range("T20").Value = 0.01291516726417287

$ws.Range("G21").Value = 1.0662095
$ws.Range("H21").Value = 2.132419
$ws.Range("I21").Value = 0.138166322546203
$ws.Range("J21").Value = 0.0965578982697073
$ws.Range("M21").Value = 24.297748
$ws.Range("N21").Value = 72.893244
$ws.Range("O21").Value = 0.1742906287498262
$ws.Range("P21").Value = 0.1859468094622229
$ws.Range("Q21").Value = 25.906489746206
$ws.Range("R21").Value = 155.438938477236
$ws.Range("S21").Value = 0.02408109522862901
$ws.Range("T21").Value = 0.01795463311162997

$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.1963773333333333
$ws.Range("H22").Value = 0.589132
$ws.Range("I22").Value = 0.02544784489173709
$ws.Range("J22").Value = 0.02667644010085692
$ws.Range("M22").Value = 30.46625333333334
$ws.Range("N22").Value = 91.39876000000001
$ws.Range("O22").Value = 0.2185380492512374
$ws.Range("P22").Value = 0.2331534018544084
$ws.Range("Q22").Value = 5.982881586257778
$ws.Range("R22").Value = 53.84593427632001
$ws.Range("S22").Value = 0.005561322380288291
$ws.Range("T22").Value = 0.006219702758880149

$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.1963773333333333
$ws.Range("H23").Value = 0.589132
$ws.Range("I23").Value = 0.02544784489173709
$ws.Range("J23").Value = 0.02667644010085692
$ws.Range("O23").Value = 0.2491807703757967
$ws.Range("P23").Value = 0.2658454419670822
$ws.Range("Q23").Value = 6.821782512650666
$ws.Range("R23").Value = 61.396042613856
$ws.Range("S23").Value = 0.006341113594526829
$ws.Range("T23").Value = 0.007091810008720702

$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.1963773333333333
$ws.Range("H24").Value = 0.589132
$ws.Range("I24").Value = 0.02544784489173709
$ws.Range("J24").Value = 0.02667644010085692
$ws.Range("M24").Value = 23.69037333333334
$ws.Range("N24").Value = 71.07112000000001
$ws.Range("O24").Value = 0.1699338582153697
$ws.Range("P24").Value = 0.181298667526812
$ws.Range("Q24").Value = 4.652252340871112
$ws.Range("R24").Value = 41.87027106784
$ws.Range("S24").Value = 0.004324450465719171
$ws.Range("T24").Value = 0.004836403044644173

$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.1963773333333333
$ws.Range("H25").Value = 0.589132
$ws.Range("I25").Value = 0.02544784489173709
$ws.Range("J25").Value = 0.02667644010085692
$ws.Range("M25").Value = 26.2168665
$ws.Range("N25").Value = 52.433733
$ws.Range("O25").Value = 0.18805669340777
$ws.Range("P25").Value = 0.1337556791894743
$ws.Range("Q25").Value = 5.148398331626
$ws.Range("R25").Value = 30.890389989756
$ws.Range("S25").Value = 0.004785637564693887
$ws.Range("T25").Value = 0.003568125364047446

$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.1963773333333333
$ws.Range("H26").Value = 0.589132
$ws.Range("I26").Value = 0.02544784489173709
$ws.Range("J26").Value = 0.02667644010085692
$ws.Range("M26").Value = 24.297748
$ws.Range("N26").Value = 72.893244
$ws.Range("O26").Value = 0.1742906287498262
$ws.Range("P26").Value = 0.1859468094622229
$ws.Range("Q26").Value = 4.771526958245333
$ws.Range("R26").Value = 42.943742624208
$ws.Range("S26").Value = 0.00443532088650891
$ws.Range("T26").Value = 0.004960398924564445
